# Update Name of Algo
# Applies updated RandomForest-imputed numeric values to specific cells
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -10.63479999999999
$ws.Range("D3").Value = -7.012999999999992
$ws.Range("E8").Value = 16.07970000000001
$ws.Range("E11").Value = 16.74599999999999
$ws.Range("A12").Value = -21.53039999999999
$ws.Range("C14").Value = -13.83079999999999
$ws.Range("E14").Value = 16.5531
$ws.Range("E15").Value = 16.07930000000001
$ws.Range("C26").Value = -12.6886
$ws.Range("D30").Value = -7.612099999999998
$ws.Range("C31").Value = -12.9672
$ws.Range("A32").Value = -21.18809999999998
$ws.Range("C35").Value = -13.20580000000001
$ws.Range("A36").Value = -19.6256
$ws.Range("E36").Value = 16.28430000000001
$ws.Range("C37").Value = -13.8135
$ws.Range("A38").Value = -19.2414
$ws.Range("D44").Value = -7.391900000000001
$ws.Range("C45").Value = -14.0049
$ws.Range("A46").Value = -21.5607
$ws.Range("A54").Value = -21.9346
$ws.Range("A55").Value = -22.45500000000001
$ws.Range("C57").Value = -14.40799999999999
$ws.Range("D58").Value = -8.536999999999994
$ws.Range("E64").Value = 17.46339999999999
$ws.Range("A67").Value = -21.53999999999998
$ws.Range("A69").Value = -21.64849999999997
$ws.Range("A72").Value = -21.64899999999998
$ws.Range("D84").Value = -8.607200000000004
$ws.Range("D89").Value = -6.193999999999996
$ws.Range("E89").Value = 18.42980000000001
$ws.Range("A91").Value = -21.5127
$ws.Range("D91").Value = -6.140199999999994
$ws.Range("D92").Value = -6.033999999999997
$ws.Range("A99").Value = -20.37259999999999
$ws.Range("C100").Value = -12.8134
$ws.Range("C102").Value = -13.92310000000001
$ws.Range("D102").Value = -7.7399
